$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(1.071259418298716, 0.8527536794795421, 1.0596515370491082, 2.832002970724148)"
$ws.Range("C2").Value = "NIG(1.5154813741650255, 1.140267621814706, 4.186618324053942, 6.301943055581725)"
$ws.Range("D2").Value = "NCT(3.0469882318656794, 1.4869425054088128, -0.3788296077850717, 2.5713244311430503)"
$ws.Range("E2").Value = "NIG(1.4070847904933308, 1.1391747949350637, 3.3961181920734598, 5.736977077035787)"
